{"js": "// Replace the date heading and the 25 two-digit multiplication prompts\n// with the updated values from the new day's worksheet.\nconst replacements = [\n  [\"2024-09-29 Sunday\", \"2024-09-30 Monday\"],\n  [\"96\u00d766=\", \"81\u00d785=\"],\n  [\"20\u00d734=\", \"85\u00d795=\"],\n  [\"71\u00d775=\", \"90\u00d789=\"],\n  [\"44\u00d763=\", \"35\u00d799=\"],\n  [\"32\u00d780=\", \"40\u00d779=\"],\n  [\"14\u00d776=\", \"96\u00d788=\"],\n  [\"20\u00d723=\", \"51\u00d773=\"],\n  [\"13\u00d760=\", \"77\u00d765=\"],\n  [\"94\u00d722=\", \"64\u00d745=\"],\n  [\"71\u00d745=\", \"67\u00d753=\"],\n  [\"56\u00d757=\", \"12\u00d783=\"],\n  [\"11\u00d787=\", \"36\u00d780=\"],\n  [\"91\u00d761=\", \"42\u00d797=\"],\n  [\"98\u00d754=\", \"32\u00d742=\"],\n  [\"71\u00d741=\", \"58\u00d729=\"],\n  [\"22\u00d798=\", \"61\u00d715=\"],\n  [\"86\u00d771=\", \"35\u00d774=\"],\n  [\"90\u00d776=\", \"92\u00d757=\"],\n  [\"52\u00d711=\", \"76\u00d720=\"],\n  [\"25\u00d725=\", \"34\u00d781=\"],\n  [\"44\u00d740=\", \"87\u00d764=\"],\n  [\"68\u00d723=\", \"51\u00d758=\"],\n  [\"83\u00d731=\", \"48\u00d796=\"],\n  [\"82\u00d725=\", \"35\u00d723=\"],\n  [\"40\u00d790=\", \"59\u00d721=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date heading and the 25 two-digit multiplication\n# prompts to the next day's generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-29 Sunday\", \"2024-09-30 Monday\"),\n    @(\"96\u00d766=\", \"81\u00d785=\"),\n    @(\"20\u00d734=\", \"85\u00d795=\"),\n    @(\"71\u00d775=\", \"90\u00d789=\"),\n    @(\"44\u00d763=\", \"35\u00d799=\"),\n    @(\"32\u00d780=\", \"40\u00d779=\"),\n    @(\"14\u00d776=\", \"96\u00d788=\"),\n    @(\"20\u00d723=\", \"51\u00d773=\"),\n    @(\"13\u00d760=\", \"77\u00d765=\"),\n    @(\"94\u00d722=\", \"64\u00d745=\"),\n    @(\"71\u00d745=\", \"67\u00d753=\"),\n    @(\"56\u00d757=\", \"12\u00d783=\"),\n    @(\"11\u00d787=\", \"36\u00d780=\"),\n    @(\"91\u00d761=\", \"42\u00d797=\"),\n    @(\"98\u00d754=\", \"32\u00d742=\"),\n    @(\"71\u00d741=\", \"58\u00d729=\"),\n    @(\"22\u00d798=\", \"61\u00d715=\"),\n    @(\"86\u00d771=\", \"35\u00d774=\"),\n    @(\"90\u00d776=\", \"92\u00d757=\"),\n    @(\"52\u00d711=\", \"76\u00d720=\"),\n    @(\"25\u00d725=\", \"34\u00d781=\"),\n    @(\"44\u00d740=\", \"87\u00d764=\"),\n    @(\"68\u00d723=\", \"51\u00d758=\"),\n    @(\"83\u00d731=\", \"48\u00d796=\"),\n    @(\"82\u00d725=\", \"35\u00d723=\"),\n    @(\"40\u00d790=\", \"59\u00d721=\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
